# Update the Saudi Arabia MSME summary sheet with refreshed (more precise)
# values for "Enterprises density (per 1000 people)" and
# "Employment (% of total)" rows, as published in the autogenerated refresh.
#
# These cells store their numbers as text (shared strings), so a plain
# assignment of a numeric-looking string (e.g. "24.85") would make Excel
# reinterpret it as a number. To keep the cell text-typed we briefly switch
# the cell to a text number format before assigning the value, then restore
# the original "Normal" style/number format so the cell's formatting is left
# exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)

    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
Set-TextValue $ws.Range("B11") "24.85"
Set-TextValue $ws.Range("C11") "3.96"
Set-TextValue $ws.Range("D11") "28.81"

# Employment (% of total): Micro (shares the same underlying shared string
# as "Enterprises density (per 1000 people)" / SMEs in the source workbook)
Set-TextValue $ws.Range("B12") "3.96"

# Employment (% of total): SMEs / MSMEs
Set-TextValue $ws.Range("C12") "15.78"
Set-TextValue $ws.Range("D12") "19.74"

Write-Output "Updated enterprise density and employment percentage figures."
